$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2: Mercury Rising | Quicksilver
$ws.Range("H2").Value = 99.55556
$ws.Range("I2").Value = 99.42856999999999
$ws.Range("K2").Value = 99.42856999999999
$ws.Range("M2").Value = 13.57143000000001

# row 5: Met a Sticky End | Animal Glue
$ws.Range("H5").Value = 170.57143
$ws.Range("I5").Value = 108.8
$ws.Range("K5").Value = 108.8
$ws.Range("M5").Value = 6.200000000000003

# row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 9124.875
$ws.Range("J62").Value = 3999.5
$ws.Range("L62").Value = 3999.5
$ws.Range("N62").Value = -5247.5

# row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 9124.875
$ws.Range("J65").Value = 3999.5
$ws.Range("L65").Value = 19997.5
$ws.Range("N65").Value = -26237.5

# row 107: Another Man's Ink | Enchanted Truegold Ink
$ws.Range("H107").Value = 3267.3
$ws.Range("I107").Value = 1963.6666
$ws.Range("K107").Value = 1963.6666
$ws.Range("M107").Value = -43.66660000000002

# row 131: Mindful Study | Grade 5 Tincture of Mind
$ws.Range("H131").Value = 9365.75
$ws.Range("I131").Value = 1856
$ws.Range("J131").Value = 19879.4
$ws.Range("K131").Value = 5568
$ws.Range("L131").Value = 59638.2
$ws.Range("M131").Value = -528
$ws.Range("N131").Value = -69718.20000000001

# row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4757.1787
$ws.Range("I141").Value = 1675.35
$ws.Range("K141").Value = 5026.049999999999
$ws.Range("M141").Value = 153.9500000000007

$ws = $wb.Worksheets.Item("ARM")
# row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 181853.5
$ws.Range("I32").Value = 254387.58
$ws.Range("K32").Value = 254387.58
$ws.Range("M32").Value = -254100.58

# row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1832.0834
$ws.Range("I45").Value = 1664.1666
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1664.1666
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1287.1666
$ws.Range("N45").Value = -2754

# row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 2328929.2
$ws.Range("I61").Value = 3337.9443
$ws.Range("K61").Value = 3337.9443
$ws.Range("M61").Value = -3125.9443

# row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1395990.5
$ws.Range("I74").Value = 1738841.4
$ws.Range("J74").Value = 24587.125
$ws.Range("K74").Value = 1738841.4
$ws.Range("L74").Value = 24587.125
$ws.Range("M74").Value = -1737967.4
$ws.Range("N74").Value = -26335.125

# row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1395990.5
$ws.Range("I77").Value = 1738841.4
$ws.Range("J77").Value = 24587.125
$ws.Range("K77").Value = 8694207
$ws.Range("L77").Value = 122935.625
$ws.Range("M77").Value = -8689839
$ws.Range("N77").Value = -131671.625

# row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 4614.1665
$ws.Range("J88").Value = 5500.125
$ws.Range("L88").Value = 5500.125
$ws.Range("N88").Value = -6312.125

# row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 4614.1665
$ws.Range("J91").Value = 5500.125
$ws.Range("L91").Value = 5500.125
$ws.Range("N91").Value = -8308.125

# row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 2077.1
$ws.Range("I122").Value = 1697.6666
$ws.Range("K122").Value = 5092.9998
$ws.Range("M122").Value = -2642.9998

# row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 664593.9
$ws.Range("I132").Value = 813938.4
$ws.Range("J132").Value = 3211.1428
$ws.Range("K132").Value = 2441815.2
$ws.Range("L132").Value = 9633.428400000001
$ws.Range("M132").Value = -2439285.2
$ws.Range("N132").Value = -14693.4284

# row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2328929.2
$ws.Range("I136").Value = 3337.9443
$ws.Range("K136").Value = 10013.8329
$ws.Range("M136").Value = -7463.832900000001

$ws = $wb.Worksheets.Item("BSM")
# row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 1424.9546
$ws.Range("J86").Value = 1179.8
$ws.Range("L86").Value = 1179.8
$ws.Range("N86").Value = -3425.8

# row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 1424.9546
$ws.Range("J89").Value = 1179.8
$ws.Range("L89").Value = 5899
$ws.Range("N89").Value = -17131

# row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 5895.971
$ws.Range("I105").Value = 5157.5835
$ws.Range("K105").Value = 5157.5835
$ws.Range("M105").Value = -3410.5835

# row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2693969.2
$ws.Range("I134").Value = 3957.1396
$ws.Range("J134").Value = 8781892
$ws.Range("K134").Value = 11871.4188
$ws.Range("L134").Value = 26345676
$ws.Range("M134").Value = -9336.418799999999
$ws.Range("N134").Value = -26350746

$ws = $wb.Worksheets.Item("CRP")
# row 44: Stay on Target | Yarzonshell Harpoon
$ws.Range("H44").Value = 4971
$ws.Range("J44").Value = 4971
$ws.Range("L44").Value = 4971
$ws.Range("N44").Value = -5855

# row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2785317.5
$ws.Range("I58").Value = 3439.5557
$ws.Range("K58").Value = 3439.5557
$ws.Range("M58").Value = -3236.5557

# row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 6090.65
$ws.Range("I105").Value = 6901.4375
$ws.Range("J105").Value = 2847.5
$ws.Range("K105").Value = 6901.4375
$ws.Range("L105").Value = 2847.5
$ws.Range("M105").Value = -5154.4375
$ws.Range("N105").Value = -6341.5

# row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 11843.333
$ws.Range("I122").Value = 2620.0715
$ws.Range("K122").Value = 7860.2145
$ws.Range("M122").Value = -5410.2145

# row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 1889.5333
$ws.Range("I132").Value = 1338.4445
$ws.Range("J132").Value = 2716.1667
$ws.Range("K132").Value = 4015.3335
$ws.Range("L132").Value = 8148.500100000001
$ws.Range("M132").Value = -1485.3335
$ws.Range("N132").Value = -13208.5001

# row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2785317.5
$ws.Range("I136").Value = 3439.5557
$ws.Range("K136").Value = 10318.6671
$ws.Range("M136").Value = -7768.667099999999

$ws = $wb.Worksheets.Item("CUL")
# row 2: Pork Is a Salty Food | Table Salt
$ws.Range("H2").Value = 148.10527
$ws.Range("I2").Value = 16.6
$ws.Range("J2").Value = 195.07143
$ws.Range("K2").Value = 99.60000000000001
$ws.Range("L2").Value = 1170.42858
$ws.Range("M2").Value = 13.39999999999999
$ws.Range("N2").Value = -1396.42858

$ws = $wb.Worksheets.Item("GSM")
# row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 5000253
$ws.Range("J2").Value = 172.75
$ws.Range("L2").Value = 172.75
$ws.Range("N2").Value = -398.75

# row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 60286.855
$ws.Range("I70").Value = 44399.8
$ws.Range("J70").Value = 100004.5
$ws.Range("K70").Value = 44399.8
$ws.Range("L70").Value = 100004.5
$ws.Range("M70").Value = -44129.8
$ws.Range("N70").Value = -100544.5

# row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 60286.855
$ws.Range("I73").Value = 44399.8
$ws.Range("J73").Value = 100004.5
$ws.Range("K73").Value = 44399.8
$ws.Range("L73").Value = 100004.5
$ws.Range("M73").Value = -43463.8
$ws.Range("N73").Value = -101876.5

# row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 3310.3076
$ws.Range("I102").Value = 3219.5
$ws.Range("K102").Value = 3219.5
$ws.Range("M102").Value = -1597.5

# row 123: Workplace Workout | Ametrine Ring of Fending
$ws.Range("H123").Value = 45860.855
$ws.Range("J123").Value = 45860.855
$ws.Range("L123").Value = 45860.855
$ws.Range("N123").Value = -50760.855

# row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 9911.571
$ws.Range("I132").Value = 7796.8335
$ws.Range("J132").Value = 22600
$ws.Range("K132").Value = 23390.5005
$ws.Range("L132").Value = 67800
$ws.Range("M132").Value = -20860.5005
$ws.Range("N132").Value = -72860

$ws = $wb.Worksheets.Item("LTW")
# row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 843.5
$ws.Range("I16").Value = 821.86664
$ws.Range("J16").Value = 908.4
$ws.Range("K16").Value = 821.86664
$ws.Range("L16").Value = 908.4
$ws.Range("M16").Value = -651.86664
$ws.Range("N16").Value = -1248.4

# row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 3503
$ws.Range("J46").Value = 5164.8
$ws.Range("L46").Value = 5164.8
$ws.Range("N46").Value = -5540.8

# row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 2124.4285
$ws.Range("I93").Value = 1724.3334
$ws.Range("J93").Value = 3124.6667
$ws.Range("K93").Value = 1724.3334
$ws.Range("L93").Value = 3124.6667
$ws.Range("M93").Value = -476.3334
$ws.Range("N93").Value = -5620.6667

# row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 3747
$ws.Range("I122").Value = 3292.2727
$ws.Range("J122").Value = 4997.5
$ws.Range("K122").Value = 9876.8181
$ws.Range("L122").Value = 14992.5
$ws.Range("M122").Value = -7426.8181
$ws.Range("N122").Value = -19892.5

# row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 11687987
$ws.Range("I132").Value = 19477618
$ws.Range("K132").Value = 58432854
$ws.Range("M132").Value = -58430324

$ws = $wb.Worksheets.Item("WVR")
# row 62: Pride Up in Smoke | Rainbow Cloth
$ws.Range("H62").Value = 9537.857
$ws.Range("I62").Value = 4929.6665
$ws.Range("J62").Value = 12994
$ws.Range("K62").Value = 4929.6665
$ws.Range("L62").Value = 12994
$ws.Range("M62").Value = -4305.6665
$ws.Range("N62").Value = -14242

# row 65: Desperate for Diversionaries (L) | Rainbow Cloth
$ws.Range("H65").Value = 9537.857
$ws.Range("I65").Value = 4929.6665
$ws.Range("J65").Value = 12994
$ws.Range("K65").Value = 24648.3325
$ws.Range("L65").Value = 64970
$ws.Range("M65").Value = -21528.3325
$ws.Range("N65").Value = -71210

# row 74: Clothing the Naked Truth | Ramie Robe of Casting
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# row 77: When in Robes (L) | Ramie Robe of Casting
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# row 100: Of Great Import | Kudzu Thread
$ws.Range("H100").Value = 4283.8184
$ws.Range("I100").Value = 3339.375
$ws.Range("K100").Value = 6678.75
$ws.Range("M100").Value = -6137.75

# row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 83335490
$ws.Range("I132").Value = 166666670
$ws.Range("K132").Value = 500000010
$ws.Range("M132").Value = -499997480

# row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 13527925
$ws.Range("I136").Value = 3106839.5
$ws.Range("J136").Value = 50001724
$ws.Range("K136").Value = 9320518.5
$ws.Range("L136").Value = 150005172
$ws.Range("M136").Value = -9317968.5
$ws.Range("N136").Value = -150010272

